# Apply updated crypto price/volume data (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "66.307.44"
$ws.Range("E2").Value = "  +2.09%  "

# Row 3
$ws.Range("D3").Value = "3.198.58"
$ws.Range("E3").Value = "  +1.45%  "

# Row 4
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.07%  "

# Row 5
$ws.Range("D5").Value = "'596.77"
$ws.Range("E5").Value = "  +3.94%  "

# Row 6
$ws.Range("D6").Value = "'154.24"
$ws.Range("E6").Value = "  +3.24%  "

# Row 7
$ws.Range("E7").Value = "  +0.03%  "

# Row 8
$ws.Range("D8").Value = "3.195.02"
$ws.Range("E8").Value = "  +1.34%  "

# Row 9
$ws.Range("D9").Value = "'0.535"
$ws.Range("E9").Value = "  +1.79%  "

# Row 10
$ws.Range("D10").Value = "'0.161"
$ws.Range("E10").Value = "  +0.63%  "

# Row 11
$ws.Range("E11").Value = "  -0.31%  "

# Row 12
$ws.Range("D12").Value = "'0.513"
$ws.Range("E12").Value = "  +3.28%  "

# Row 14
$ws.Range("D14").Value = "'38.97"
$ws.Range("E14").Value = "  +5.23%  "

# Row 15
$ws.Range("D15").Value = "3.724.78"
$ws.Range("E15").Value = "  +1.52%  "

# Row 16
$ws.Range("D16").Value = "66.263.26"
$ws.Range("E16").Value = "  +1.86%  "

# Row 17
$ws.Range("D17").Value = "'7.44"
$ws.Range("E17").Value = "  +4.85%  "

# Row 18
$ws.Range("D18").Value = "3.207.77"
$ws.Range("E18").Value = "  +1.70%  "

# Row 19
$ws.Range("D19").Value = "'0.112"
$ws.Range("E19").Value = "  +0.62%  "

# Row 20
$ws.Range("D20").Value = "'510.27"
$ws.Range("E20").Value = "  +0.84%  "

# Row 21
$ws.Range("D21").Value = "'15.34"
$ws.Range("E21").Value = "  +3.88%  "

# Row 22
$ws.Range("D22").Value = "'0.737"
$ws.Range("E22").Value = "  +2.99%  "

# Row 23
$ws.Range("D23").Value = "'8.00"
$ws.Range("E23").Value = "  +3.55%  "

# Row 24
$ws.Range("D24").Value = "'15.12"
$ws.Range("E24").Value = "  -1.12%  "

# Row 25
$ws.Range("E25").Value = "  +0.83%  "

# Row 26
$ws.Range("E26").Value = "  -0.01%  "

# Row 27
$ws.Range("D27").Value = "'9.30"
$ws.Range("E27").Value = "  +5.38%  "

# Row 28
$ws.Range("E28").Value = "  +3.10%  "

# Row 29
$ws.Range("E29").Value = "  +5.66%  "

# Row 30
$ws.Range("D30").Value = "'2.92"
$ws.Range("E30").Value = "  +4.31%  "

# Row 31
$ws.Range("D31").Value = "'6.96"
$ws.Range("E31").Value = "  +12.46%  "

# Row 32
$ws.Range("D32").Value = "'28.23"
$ws.Range("E32").Value = "  +2.33%  "

# Row 33
$ws.Range("E33").Value = "  +3.29%  "

# Row 34
$ws.Range("E34").Value = "  +0.03%  "

# Row 35
$ws.Range("D35").Value = "'6.53"
$ws.Range("E35").Value = "  +0.33%  "

# Row 36
$ws.Range("D36").Value = "'54.78"
$ws.Range("E36").Value = "  -0.30%  "

# Row 37
$ws.Range("D37").Value = "'0.0900"
$ws.Range("E37").Value = "  +0.30%  "

# Row 38
$ws.Range("D38").Value = "'484.94"
$ws.Range("E38").Value = "  +4.42%  "

# Row 39
$ws.Range("D39").Value = "'0.0418"
$ws.Range("E39").Value = "  -0.41%  "

# Row 40
$ws.Range("B40").Value = "Cosmos"
$ws.Range("C40").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D40").Value = "'8.84"
$ws.Range("E40").Value = "  +2.20%  "

# Row 41
$ws.Range("B41").Value = "dogwifhat"
$ws.Range("C41").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D41").Value = "'2.91"
$ws.Range("E41").Value = "  -2.48%  "

# Row 42
$ws.Range("D42").Value = "'0.122"
$ws.Range("E42").Value = "  +4.93%  "

# Row 43
$ws.Range("D43").Value = "'0.299"
$ws.Range("E43").Value = "  +6.27%  "

# Row 44
$ws.Range("E44").Value = "  +12.40%  "

# Row 45
$ws.Range("D45").Value = "2.924.74"
$ws.Range("E45").Value = "  -4.10%  "

# Row 46
$ws.Range("D46").Value = "'2.41"
$ws.Range("E46").Value = "  -0.64%  "

# Row 47
$ws.Range("D47").Value = "'28.47"
$ws.Range("E47").Value = "  -0.16%  "

# Row 48
$ws.Range("E48").Value = "  -0.04%  "

# Row 49
$ws.Range("E49").Value = "  +2.04%  "

# Row 50
$ws.Range("D50").Value = "'2.32"
$ws.Range("E50").Value = "  +3.17%  "

# Row 51
$ws.Range("D51").Value = "'2.60"
$ws.Range("E51").Value = "  +5.39%  "

